# Fixes until github call
# Insert two new cluster rows ("conduction" and "zaakonline") above the
# existing "testcluster" row on the "clusters" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("clusters")

# Push the existing row (testcluster, ...) down by two rows, making room
# for the two new rows at the top.
$ws.Rows("1:2").Insert()

# Row 1: conduction
$ws.Range("A1").Value = "conduction"
$ws.Range("B1").Value = "the main conduction cluster"
$ws.Range("C1").Value = "the main conduction domain"
$ws.Range("D1").Value = "conduction.nl"

# Row 2: zaakonline
$ws.Range("A2").Value = "zaakonline"
$ws.Range("B2").Value = "the main zaakonline cluster"
$ws.Range("C2").Value = "the main zaakonline domain"
$ws.Range("D2").Value = "zaakonline.nl"

# Match the saved selection state (A1:D2 highlighted, active cell on row 2).
$ws.Activate()
$ws.Range("A1:D2").Select()
